$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Device key" row content (shared string no longer referenced in v2 of the lab)
$ws.Range("B10").ClearContents()
